# PurchaseTable.xlsx - "Inventory, Reward 구현"
#
# Inserts a new "rewardId" column right after "internalProductId" (new
# column B), pushing the existing kind/title/isActive/storeSkuApple/
# storeSkuGoogle block one column to the right (old B:F -> new C:G), and
# fills in the new column's header, type row and description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot every existing value *before* anything is overwritten, since
# the shift is going to write over the very cells we still need to read.
$old = @{}
foreach ($r in 1..6) {
  foreach ($col in @("B","C","D","E","F")) {
    $old["$col$r"] = $ws.Range("$col$r").Value()
  }
}

# Shift old column B..F one column right, into new C..G.
foreach ($r in 1..6) {
  $ws.Range("C$r").Value = $old["B$r"]
  $ws.Range("D$r").Value = $old["C$r"]
  $ws.Range("E$r").Value = $old["D$r"]
  $ws.Range("F$r").Value = $old["E$r"]
  $ws.Range("G$r").Value = $old["F$r"]
}

# Populate the new "rewardId" column (B).
$ws.Range("B1").Value = "rewardId"
$ws.Range("B2").Value = "string"
$ws.Range("B3").Value = $null
$ws.Range("B4").Value = "Reward Key (REWARD.rewardId(pk))"
$ws.Range("B5").Value = $null
$ws.Range("B6").Value = $null

# Match the new last column's display width to the (former) last column's.
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth
